$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(4.612846683653662, 9.413972823782968, 9.413972823782968)
    3  = @(9.981690992661678, 18.48461294937369, 18.48461294937369)
    4  = @(10.21802379466729, 16.48068353978608, 16.48068353978608)
    5  = @(10.18162570406218, 16.16131064136858, 16.16131064136858)
    6  = @(10.29058658973594, 14.70083798533697, 14.70083798533697)
    7  = @(10.09243780125401, 14.41776828750572, 14.41776828750572)
    8  = @(10.0300599756187,  13.73980818577903, 13.73980818577903)
    9  = @(10.10764763466389, 13.84609265022454, 13.84609265022454)
    10 = @(10.10819953617945, 14.03916602247153, 14.03916602247153)
    11 = @(10.04802561743162, 13.95559113532174, 13.95559113532174)
    12 = @(9.964361748524935, 13.64981061441747, 13.64981061441747)
    13 = @(9.935132616278832, 13.60977070723142, 13.60977070723142)
    14 = @(9.877620505261458, 13.53098699350885, 13.53098699350885)
    15 = @(9.730356576858799, 14.96977934901378, 14.96977934901378)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}
